# "Fix test case ID" - the form's title/id strings on the "settings" sheet
# were stamped with a stale date (20230209); update them to 20210331.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

$ws.Range("A2").Value = "Household Exit Survey – 20210331"
$ws.Range("B2").Value = "case_deactivate_20210331"

# Mirror the author's click into B2 (the cell they just retyped) so the
# saved selection matches.
$ws.Range("B2").Select() | Out-Null
